$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Saturday hours for the week of 2018-04-16 (row 13) from 6 to 8.
$ws.Range("B13").Value = 8

# Move the active cell selection to E14, matching where editing continued.
$ws.Range("E14").Select()
